# Inserts a new weekly price record as row 163 in the "Zapallo italiano"
# sheet, pushing the existing rows 163:190 down to 164:191 (dimension
# grows from A1:R190 to A1:R191).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 163; Excel shifts rows 163:190 -> 164:191.
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new record's data.
$ws.Cells.Item(163, 1).Value2 = 11
$ws.Cells.Item(163, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(163, 3).Value2 = "Bíobío"
$ws.Cells.Item(163, 4).Value2 = 44964
$ws.Cells.Item(163, 5).Value2 = 8
$ws.Cells.Item(163, 6).Value2 = 100112032
$ws.Cells.Item(163, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(163, 8).Value2 = "Sin especificar"
$ws.Cells.Item(163, 9).Value2 = "Primera"
$ws.Cells.Item(163, 10).Value2 = 150
$ws.Cells.Item(163, 11).Value2 = 8500
$ws.Cells.Item(163, 12).Value2 = 9000
$ws.Cells.Item(163, 13).Value2 = 8667
$ws.Cells.Item(163, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(163, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(163, 16).Value2 = 173
$ws.Cells.Item(163, 17).Value2 = 50
$ws.Cells.Item(163, 18).Value2 = "Hortaliza"
